$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, pushing the existing "Vincent Huang" entry down to row 3.
# Excel's default row insert copies the formatting of the row above (the bold/bordered
# header row), so we need to normalize the new row's formatting afterwards.
$ws.Rows.Item(2).Insert()

# B2:D2 should use the plain/default style (no border, no bold) like the data rows.
$ws.Range("B2:D2").Style = "Normal"

# A2 (the name cell) should use the same bordered/bold/centered style as the other
# name cells in column A. Copy the formatting from A3 (the original name-cell style).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Populate the new row 2 with the "Adam Jackson" entry.
$ws.Range("A2").Value = "Adam Jackson"
$ws.Range("B2").Value = "Licensed Administrator"
$ws.Range("C2").Value = "Vancouver, BC"
$ws.Range("D2").Value = "RE/MAX Select Properties"

# Row 3 (previously row 2) keeps its title/city/office, but the name changes.
$ws.Range("A3").Value = "Sarka Trileta"
